$p = $ppt.ActivePresentation

# Locate the slide/shape containing the code comment that needs correcting.
# (The slide about array indexing — "Array elements are accessed using
# bracket notation." — currently reads "...of a1 (the first integer)"
# and should read "...of a (the first integer)".)
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t.Contains("integer at index 0 of a1")) {
                $targetSlide = $s
                $targetShape = $sh
            }
        }
    }
}

if ($targetShape -ne $null) {
    $tr = $targetShape.TextFrame.TextRange
    $full = $tr.Text

    # Find "of a1 " immediately preceding "(the first integer)" and retype it
    # as "of a " -- this mirrors selecting "of a1 " in the editor and typing
    # "of a " over it, which is the minimal correction matching the commit
    # ("minor correction to one slide").
    $needle = "of a1 "
    $idx = $full.IndexOf($needle)
    if ($idx -ge 0) {
        $start = $idx + 1   # TextRange.Characters is 1-based
        $len = $needle.Length
        $sub = $tr.Characters($start, $len)
        $sub.Text = "of a "
    }
}
